$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Lista Tipo Logradouro")

# The template sheet ends up with its whole grid selected (as it did when
# it was last used as the basis for copying the new "Lista" sheets).
$template.Cells.Select() | Out-Null

# --- New sheet 18: "Lista Classe Social" ---
$template.Copy($null, $template)
$sheet18 = $wb.Worksheets.Item($template.Index + 1)
$sheet18.Name = "Lista Classe Social"
$sheet18.Range("B1").Value = "linha tipo 18"
$sheet18.Range("C5").Value = 10
$sheet18.Cells.Select() | Out-Null

# --- New sheet 19: "Lista Tipo de Uso do Imovel" ---
$template.Copy($null, $sheet18)
$sheet19 = $wb.Worksheets.Item($sheet18.Index + 1)
$sheet19.Name = "Lista Tipo de Uso do Imovel"
$sheet19.Range("B1").Value = "linha tipo 19"
$sheet19.Range("C5").Value = 20
$sheet19.Cells.Select() | Out-Null

# --- New sheet 20: "Lista Acesso ao Hridrometro" ---
$template.Copy($null, $sheet19)
$sheet20 = $wb.Worksheets.Item($sheet19.Index + 1)
$sheet20.Name = "Lista Acesso ao Hridrometro"
$sheet20.Range("B1").Value = "linha tipo 20"
$sheet20.Range("C5").Value = 20
$sheet20.Range("B5").Select() | Out-Null

foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
